$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 8 with Quick sort problem entry
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Quick sort"
$ws.Range("D8").Value = "c"
$ws.Range("E8").Value = "O(n2)"
$ws.Range("F8").Value = "O(1)"
$ws.Range("H8").Value = "QuickSort"

# Update selection to match target state
$ws.Range("H8").Select()
